# Generate Report for Handoff
# Updates the status/handoff timestamps for file "4718eab0-2d09-4fb1-a4b5-9300d243b4b9.md"
# which moved from "In Translation" to "Ready for handoff", refreshing the
# "Latest Handoff Date(time)" columns on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "2016-03-23 10:26:44"
$ws.Range("D10").Value = "2016-03-23 10:26:44"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("E6").Value = "2016-03-23 10:26:40"
$ws.Range("E9").Value = "2016-03-23 10:26:40"
$ws.Range("E10").Value = "2016-03-23 10:26:40"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("E6").Value = "2016-03-23 10:26:44"
$ws.Range("E9").Value = "2016-03-23 10:26:44"
$ws.Range("E10").Value = "2016-03-23 10:26:44"
